$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '61.901.02'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = "'" + '3.413.82'
$ws.Range("E3").Value = '  -0.37%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").Value = "'" + '408.28'
$ws.Range("E5").Value = '  -0.06%  '

$ws.Range("D6").Value = "'" + '128.62'
$ws.Range("E6").Value = '  -3.75%  '

$ws.Range("D7").Value = "'" + '0.633'
$ws.Range("E7").Value = '  +6.61%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = "'" + '0.739'
$ws.Range("E9").Value = '  +9.00%  '

$ws.Range("E10").Value = '  +15.44%  '

$ws.Range("D11").Value = "'" + '42.53'
$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = "'" + '0.141'
$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = "'" + '3.965.56'
$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'" + '8.88'
$ws.Range("E14").Value = '  +4.92%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = "'" + '0.0000209'
$ws.Range("E15").Value = '  +57.02%  '

$ws.Range("D16").Value = "'" + '20.90'
$ws.Range("E16").Value = '  +4.55%  '

$ws.Range("D17").Value = "'" + '3.406.67'
$ws.Range("E17").Value = '  -0.71%  '

$ws.Range("D18").Value = "'" + '12.11'
$ws.Range("E18").Value = '  +9.73%  '

$ws.Range("E19").Value = '  +3.36%  '

$ws.Range("D20").Value = "'" + '61.859.93'
$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").Value = "'" + '406.85'
$ws.Range("E21").Value = '  +29.04%  '

$ws.Range("D22").Value = "'" + '89.55'
$ws.Range("E22").Value = '  +5.18%  '

$ws.Range("D24").Value = "'" + '13.07'
$ws.Range("E24").Value = '  +1.94%  '

$ws.Range("D25").Value = "'" + '3.22'
$ws.Range("E25").Value = '  +3.00%  '

$ws.Range("D26").Value = "'" + '32.69'
$ws.Range("E26").Value = '  +9.66%  '

$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = "'" + '4.80'
$ws.Range("E27").Value = '  +0.20%  '

$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = "'" + '8.57'
$ws.Range("E28").Value = '  +3.61%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'" + '7.58'
$ws.Range("E29").Value = '  -1.65%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'" + '2.74'
$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").Value = "'" + '0.119'
$ws.Range("E31").Value = '  +2.08%  '

$ws.Range("D32").Value = "'" + '0.172'
$ws.Range("E32").Value = '  -1.54%  '

$ws.Range("D33").Value = "'" + '11.83'
$ws.Range("E33").Value = '  +3.49%  '

$ws.Range("D34").Value = "'" + '43.00'
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("E35").Value = '  +0.67%  '

$ws.Range("D36").Value = "'" + '0.0495'
$ws.Range("E36").Value = '  +1.84%  '

$ws.Range("D37").Value = "'" + '54.10'
$ws.Range("E37").Value = '  +3.76%  '

$ws.Range("D38").Value = "'" + '0.998'
$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").Value = "'" + '3.36'
$ws.Range("E39").Value = '  -2.39%  '

$ws.Range("E40").Value = '  +6.92%  '

$ws.Range("E41").Value = '  -3.08%  '

$ws.Range("E42").Value = '  +5.91%  '

$ws.Range("D43").Value = "'" + '140.87'
$ws.Range("E43").Value = '  +2.05%  '

$ws.Range("E44").Value = '  -2.14%  '

$ws.Range("D45").Value = "'" + '4.05'
$ws.Range("E45").Value = '  +1.20%  '

$ws.Range("E46").Value = '  +8.63%  '

$ws.Range("D47").Value = "'" + '16.54'
$ws.Range("E47").Value = '  -2.13%  '

$ws.Range("D48").Value = "'" + '21.81'
$ws.Range("E48").Value = '  +1.41%  '

$ws.Range("D49").Value = "'" + '2.119.29'
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").Value = "'" + '2.39'
$ws.Range("E50").Value = '  +4.12%  '

$ws.Range("D51").Value = "'" + '0.133'
$ws.Range("E51").Value = '  +17.82%  '

